# Update the "dSF" column (F) values for the davies_zach 2021 save-data sheet.
# This mirrors a re-pull of the upstream data followed by pushing the refreshed
# values back into the sheet (per commit message: "repull data, push all data,
# mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 2
    4  = -4
    5  = 2
    8  = -2
    9  = -2
    10 = -7
    11 = -5
    12 = -4
    13 = 7
    14 = -2
    15 = -2
    17 = -1
    18 = 3
    19 = -6
    20 = 2
    21 = 5
    23 = 2
    27 = 3
    29 = -1
    30 = -2
    31 = -3
    32 = -4
    33 = 1
    34 = 2
    35 = 4
    36 = 5
    37 = 1
    38 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
